# Sanity.xlsx edit:
#   "working days november run 2 without probation, without end of year"
#
# The previous row 39 (TCID 39 - DeactivationLeaveBalance, pointing at the
# old "deactivation//Leave_Scenarios_Without_Creation.xlsx" class/file) is
# removed, which shifts rows 40-43 up to rows 39-42. A brand new row 43 is
# then appended, re-introducing a "DeactivationLeaveBalance" test case
# (TCID 44) that now points at the new Accural//Deactivation.xlsx
# class/workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 39 - this shifts rows 40:43 up to 39:42.
$null = $ws.Rows("39").Delete()

# Append the new row 43 with the refreshed Deactivation test case data.
$ws.Range("A43").Value = "44"
$ws.Range("B43").Value = "DeactivationLeaveBalance"
$ws.Range("C43").Value = "DeactivationLeaveBalance"
$ws.Range("D43").Value = "com.darwinbox.leaves.Accural.Custom.Deactivation"
$ws.Range("E43").Value = "Accural//Deactivation.xlsx"
$ws.Range("F43").Value = "All_without_Creation"
$ws.Range("G43").Value = "All"

# Match the text number-format used by the rest of column A / G.
$ws.Range("A43").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"

# Restore the view: scrolled so row 28 is at the top, with the newly
# shifted-up row 39 selected as a full row.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Rows("39").Select()
